$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
# These cells hold price strings such as "211.60" or "26.222.04". Several of
# them contain only a single "." so Excel would otherwise parse them as a
# plain number (losing the trailing zero, e.g. "211.60" -> 211.6). Force the
# cell to Text format first so the literal digits are preserved exactly.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.222.04"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.672.46"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.60"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5278"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2645"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06278"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.34"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07561"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.675.85"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.461"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5600"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.92"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008008"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.025.67"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.814"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "187.59"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.40"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.223"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.86"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.578"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.96"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06226"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.360"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.501"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.632"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.002"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6053"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.411"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.749"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.116"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01620"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.101.70"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8743"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.73"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.824.72"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000112"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.027"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4257"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.984"

# --- Volume(1h) (column E) updates ---
# These are already non-numeric-looking strings (padded with spaces / % sign)
# so a plain Value assignment keeps them as text without touching styles.
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("E6").Value = "  -3.89%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("E8").Value = "  -3.39%  "
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("E10").Value = "  -3.03%  "
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("E13").Value = "  -2.04%  "
$ws.Range("E14").Value = "  -4.44%  "
$ws.Range("E16").Value = "  -4.94%  "
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("E21").Value = "  -5.44%  "
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("E25").Value = "  -4.87%  "
$ws.Range("E26").Value = "  -4.30%  "
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("E30").Value = "  -3.59%  "
$ws.Range("E31").Value = "  -3.19%  "
$ws.Range("E32").Value = "  -4.62%  "
$ws.Range("E33").Value = "  -3.53%  "
$ws.Range("E34").Value = "  -3.73%  "
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("E45").Value = "  +3.20%  "
$ws.Range("E46").Value = "  -2.73%  "
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("E51").Value = "  -2.35%  "
